# Auto-generated edit script applying updated market price data
# to the Maduin_Profits workbook sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 999
$ws.Range("J9").Value = 999
$ws.Range("L9").Value = 999
$ws.Range("N9").Value = -1337
$ws.Range("H17").Value = 506.125
$ws.Range("J17").Value = 539.8
$ws.Range("L17").Value = 1619.4
$ws.Range("N17").Value = -1955.4
$ws.Range("H100").Value = 2427.2856
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2598.2
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2598.2
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3680.2
$ws.Range("H111").Value = 541.8182
$ws.Range("I111").Value = 286.8
$ws.Range("K111").Value = 860.4000000000001
$ws.Range("M111").Value = 2206.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 472.75
$ws.Range("I5").Value = 300.4
$ws.Range("K5").Value = 300.4
$ws.Range("M5").Value = -188.4
$ws.Range("H32").Value = 3292
$ws.Range("I32").Value = 3292
$ws.Range("K32").Value = 3292
$ws.Range("M32").Value = -3005
$ws.Range("H61").Value = 3337.2
$ws.Range("J61").Value = 6302.8
$ws.Range("L61").Value = 6302.8
$ws.Range("N61").Value = -6726.8
$ws.Range("H74").Value = 2301.2856
$ws.Range("J74").Value = 2599.75
$ws.Range("L74").Value = 2599.75
$ws.Range("N74").Value = -4347.75
$ws.Range("H77").Value = 2301.2856
$ws.Range("J77").Value = 2599.75
$ws.Range("L77").Value = 12998.75
$ws.Range("N77").Value = -21734.75
$ws.Range("H88").Value = 1722.9412
$ws.Range("I88").Value = 1442.1428
$ws.Range("J88").Value = 1919.5
$ws.Range("K88").Value = 1442.1428
$ws.Range("L88").Value = 1919.5
$ws.Range("M88").Value = -1036.1428
$ws.Range("N88").Value = -2731.5
$ws.Range("H91").Value = 1722.9412
$ws.Range("I91").Value = 1442.1428
$ws.Range("J91").Value = 1919.5
$ws.Range("K91").Value = 1442.1428
$ws.Range("L91").Value = 1919.5
$ws.Range("M91").Value = -38.14280000000008
$ws.Range("N91").Value = -4727.5
$ws.Range("H102").Value = 800
$ws.Range("I102").Value = 800
$ws.Range("K102").Value = 800
$ws.Range("M102").Value = 822
$ws.Range("H122").Value = 494.66666
$ws.Range("I122").Value = 494.66666
$ws.Range("K122").Value = 1483.99998
$ws.Range("M122").Value = 966.0000199999999
$ws.Range("H136").Value = 3337.2
$ws.Range("J136").Value = 6302.8
$ws.Range("L136").Value = 18908.4
$ws.Range("N136").Value = -24008.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 472.75
$ws.Range("I4").Value = 300.4
$ws.Range("K4").Value = 300.4
$ws.Range("M4").Value = -185.4
$ws.Range("H86").Value = 5578.2856
$ws.Range("I86").Value = 2012.25
$ws.Range("J86").Value = 10333
$ws.Range("K86").Value = 2012.25
$ws.Range("L86").Value = 10333
$ws.Range("M86").Value = -889.25
$ws.Range("N86").Value = -12579
$ws.Range("H89").Value = 5578.2856
$ws.Range("I89").Value = 2012.25
$ws.Range("J89").Value = 10333
$ws.Range("K89").Value = 10061.25
$ws.Range("L89").Value = 51665
$ws.Range("M89").Value = -4445.25
$ws.Range("N89").Value = -62897
$ws.Range("H99").Value = 2077.1
$ws.Range("I99").Value = 1733.4375
$ws.Range("K99").Value = 1733.4375
$ws.Range("M99").Value = -235.4375
$ws.Range("H105").Value = 867.25
$ws.Range("I105").Value = 867.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 867.25
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 879.75
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3633.1667
$ws.Range("I22").Value = 783
$ws.Range("J22").Value = 6483.3335
$ws.Range("K22").Value = 783
$ws.Range("L22").Value = 6483.3335
$ws.Range("M22").Value = -433
$ws.Range("N22").Value = -7183.3335
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("K76").Value = 2000
$ws.Range("M76").Value = -1685
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("K79").Value = 2000
$ws.Range("M79").Value = -908
$ws.Range("H99").Value = 2500.6667
$ws.Range("I99").Value = 3051
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 3051
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = -1553
$ws.Range("N99").Value = -4396
$ws.Range("H126").Value = 2500.6667
$ws.Range("I126").Value = 3051
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 9153
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -6683
$ws.Range("N126").Value = -9140
$ws.Range("H134").Value = 1443.6364
$ws.Range("I134").Value = 1474.6666
$ws.Range("J134").Value = 1304
$ws.Range("K134").Value = 4423.9998
$ws.Range("L134").Value = 3912
$ws.Range("M134").Value = -1888.9998
$ws.Range("N134").Value = -8982

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9000
$ws.Range("J113").Value = 9000
$ws.Range("L113").Value = 9000
$ws.Range("N113").Value = -13340
$ws.Range("H122").Value = 23443.25
$ws.Range("I122").Value = 28476.54
$ws.Range("K122").Value = 85429.62
$ws.Range("M122").Value = -82979.62
$ws.Range("H132").Value = 2530.5
$ws.Range("I132").Value = 2419.1
$ws.Range("J132").Value = 2809
$ws.Range("K132").Value = 7257.299999999999
$ws.Range("L132").Value = 8427
$ws.Range("M132").Value = -4727.299999999999
$ws.Range("N132").Value = -13487

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2309.0908
$ws.Range("I22").Value = 1625
$ws.Range("J22").Value = 2700
$ws.Range("K22").Value = 1625
$ws.Range("L22").Value = 2700
$ws.Range("M22").Value = -1330
$ws.Range("N22").Value = -3290
$ws.Range("H27").Value = 2309.0908
$ws.Range("I27").Value = 1625
$ws.Range("J27").Value = 2700
$ws.Range("K27").Value = 1625
$ws.Range("L27").Value = 2700
$ws.Range("M27").Value = -1518
$ws.Range("N27").Value = -2914
$ws.Range("H132").Value = 2249.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 465
$ws.Range("I81").Value = 465
$ws.Range("K81").Value = 930
$ws.Range("M81").Value = 131
$ws.Range("H84").Value = 465
$ws.Range("I84").Value = 465
$ws.Range("K84").Value = 4650
$ws.Range("M84").Value = 654
$ws.Range("H107").Value = 536.6667
$ws.Range("I107").Value = 284.2
$ws.Range("J107").Value = 1799
$ws.Range("K107").Value = 852.5999999999999
$ws.Range("L107").Value = 5397
$ws.Range("M107").Value = 1067.4
$ws.Range("N107").Value = -9237
$ws.Range("H132").Value = 1750.5555
$ws.Range("I132").Value = 1719.375
$ws.Range("K132").Value = 5158.125
$ws.Range("M132").Value = -2628.125
$ws.Range("H140").Value = 80000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 80000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 80000
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -90360
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
